# Data.xlsx update: "add reconnect plc feature, add API filter"
#
# The PLC reconnect logic re-reads the counter/runtime registers for Line 1
# (row 2 of Sheet1), so its logged "counterOut" (column D) and "runtime"
# (column G) values move forward. Both columns are stored as text (shared
# strings) in this workbook, so we have to land the new values as literal
# text rather than numbers - writing straight into .Value would get
# silently re-typed as a number by Excel. The reliable way to force a
# literal-text cell without picking up a new number-format style (which
# would show up as an unwanted style diff) is: stage the text via a
# TEXT() formula, then Copy / PasteSpecial-values over itself so the
# formula collapses to a plain inline literal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: counterOut 580 -> 12727
$ws.Range("D2").Formula = '=TEXT(12727,"0")'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

# G2: runtime 22.0 -> 214.0
$ws.Range("G2").Formula = '=TEXT(214,"0.0")'
$ws.Range("G2").Copy()
$ws.Range("G2").PasteSpecial(-4163)
